$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 10).Value = 8
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(5, 10).Value = 3
$ws.Cells.Item(6, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 2
$ws.Cells.Item(8, 10).Value = 3
$ws.Cells.Item(9, 10).Value = 8
$ws.Cells.Item(10, 10).Value = 4
$ws.Cells.Item(11, 10).Value = 2
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(13, 10).Value = 7
$ws.Cells.Item(14, 10).Value = 2
$ws.Cells.Item(15, 10).Value = 2
$ws.Cells.Item(16, 10).Value = 7
$ws.Cells.Item(17, 10).Value = 6
$ws.Cells.Item(18, 10).Value = 1
$ws.Cells.Item(19, 10).Value = 1
$ws.Cells.Item(21, 10).Value = 3
$ws.Cells.Item(24, 10).Value = 7
$ws.Cells.Item(25, 10).Value = 7
$ws.Cells.Item(26, 10).Value = 3
$ws.Cells.Item(27, 10).Value = 2
$ws.Cells.Item(28, 10).Value = 7
$ws.Cells.Item(29, 10).Value = 4
$ws.Cells.Item(30, 10).Value = 1
$ws.Cells.Item(31, 10).Value = 8
$ws.Cells.Item(32, 10).Value = 1
$ws.Cells.Item(33, 10).Value = 4
$ws.Cells.Item(35, 10).Value = 5
$ws.Cells.Item(37, 10).Value = 7
$ws.Cells.Item(38, 10).Value = 1
$ws.Cells.Item(40, 10).Value = 8
$ws.Cells.Item(41, 10).Value = 7
$ws.Cells.Item(42, 10).Value = 6
$ws.Cells.Item(43, 10).Value = 4
$ws.Cells.Item(44, 10).Value = 5
$ws.Cells.Item(45, 10).Value = 5
$ws.Cells.Item(46, 10).Value = 1
$ws.Cells.Item(47, 10).Value = 7
$ws.Cells.Item(48, 10).Value = 7
$ws.Cells.Item(49, 10).Value = 8
$ws.Cells.Item(51, 10).Value = 5
$ws.Cells.Item(52, 10).Value = 1
$ws.Cells.Item(53, 10).Value = 3
$ws.Cells.Item(54, 10).Value = 5
$ws.Cells.Item(57, 10).Value = 4
$ws.Cells.Item(58, 10).Value = 7
$ws.Cells.Item(59, 10).Value = 4
$ws.Cells.Item(60, 10).Value = 2
$ws.Cells.Item(63, 10).Value = 8
$ws.Cells.Item(64, 10).Value = 2
$ws.Cells.Item(65, 10).Value = 2
$ws.Cells.Item(68, 10).Value = 8
$ws.Cells.Item(69, 10).Value = 5
$ws.Cells.Item(70, 10).Value = 3
$ws.Cells.Item(71, 10).Value = 7
$ws.Cells.Item(72, 10).Value = 3
$ws.Cells.Item(73, 10).Value = 4
$ws.Cells.Item(74, 10).Value = 4
$ws.Cells.Item(75, 10).Value = 3
$ws.Cells.Item(77, 10).Value = 3
$ws.Cells.Item(78, 10).Value = 3
$ws.Cells.Item(79, 10).Value = 2
$ws.Cells.Item(80, 10).Value = 1
$ws.Cells.Item(81, 10).Value = 4
